# Highlight (yellow) the "Adversary AI" implementation sub-bullets:
#   K-means clustering algorithm / K-means expenation / Vision radius /
#   Point saving / Point loading / Clusterising / NavGoal
# i.e. every paragraph strictly between the "Adversary AI" bullet and the
# "Problems and solutions:" bullet. Both the run text and the paragraph
# mark itself (pPr/rPr) need the highlight, so we go through
# Paragraph.Range.Font rather than Paragraph.Range directly.

$d = $word.ActiveDocument
$wdYellow = 7  # wdColorIndex.wdYellow

$inBlock = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.Trim()

    if ($t -eq "Problems and solutions:") {
        break
    }

    if ($inBlock) {
        $p.Range.Font.HighlightColorIndex = $wdYellow
    }

    if ($t -eq "Adversary AI-") {
        $inBlock = $true
    }
}
